$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric/date values in row 2
$ws.Range("A2").Value = 46045
$ws.Range("B2").Value = 31.49
$ws.Range("C2").Value = 23.74
$ws.Range("D2").Value = 20.57
$ws.Range("E2").Value = 15.91
$ws.Range("F2").Value = 15.89
$ws.Range("G2").Value = 24.91
$ws.Range("H2").Value = 35.48
$ws.Range("I2").Value = 58.29
$ws.Range("J2").Value = 56.83
$ws.Range("K2").Value = 31.82
$ws.Range("L2").Value = 6.61
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 0.91
$ws.Range("O2").Value = 0.92
$ws.Range("P2").Value = 2.02
$ws.Range("Q2").Value = 4.25
$ws.Range("R2").Value = 14.47
$ws.Range("S2").Value = 43.94
$ws.Range("T2").Value = 80.64
$ws.Range("U2").Value = 101
$ws.Range("V2").Value = 107.1
$ws.Range("W2").Value = 100.25
$ws.Range("X2").Value = 81.83
$ws.Range("Y2").Value = 61.61
$ws.Range("Z2").Value = 38.4

# Slot 4h max/price
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 87.7

# Slot 2h first/price
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 103.68

# Slot 2h second/price
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 90.82

# Slot min price
$ws.Range("AG2").Value = "0h-16h"
